# Rename tabs in Eldata.xlsx so that elasticity/productivity tabs are
# named consistently per production structure (KL), and make the
# renamed "prodKL" tab (previously "FPROD") the active tab.

$wb = $excel.ActiveWorkbook

$wsElasPROD = $wb.Worksheets.Item("elasPROD")
$wsElasPROD.Name = "elasKL"

$wsFPROD = $wb.Worksheets.Item("FPROD")
$wsFPROD.Name = "prodKL"

# Make the renamed "prodKL" sheet (4th tab) the active / selected tab,
# replacing "elasFU" as the previously selected tab.
$wsFPROD.Activate()
